$wb = $excel.ActiveWorkbook

# Rename sheets (truncate trailing characters as in the target diff)
$wb.Worksheets.Item("Include from MeasurePopulationT").Name = "Include from MeasurePopulatio"
$wb.Worksheets.Item("Include from Measure Population").Name = "Include from Measure Populati"
$wb.Worksheets.Item("Exclude from MeasurePopulationT").Name = "Exclude from MeasurePopulatio"

# Update the Date value on the Metadata sheet
$ws = $wb.Worksheets.Item("Metadata")
$ws.Range("B8").Value = "2021-10-01T15:07:10+00:00"
